$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '19.729.81'
$ws.Cells.Item(2, 5).Value = '  -8.91%  '
$ws.Cells.Item(3, 4).Value = '1.385.58'
$ws.Cells.Item(3, 5).Value = '  -9.71%  '
$c = $ws.Cells.Item(4, 4)
$c.Value = "'1.004"
$c.Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.39%  '
$c = $ws.Cells.Item(5, 4)
$c.Value = "'1.003"
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.27%  '
$c = $ws.Cells.Item(6, 4)
$c.Value = "'267.66"
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -7.20%  '
$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.3634"
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -7.37%  '
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.3031"
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -4.42%  '
$c = $ws.Cells.Item(9, 4)
$c.Value = "'37.77"
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -11.05%  '
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.9709"
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -7.73%  '
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.06396"
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -10.87%  '
$c = $ws.Cells.Item(12, 4)
$c.Value = "'1.004"
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +0.37%  '
$c = $ws.Cells.Item(13, 4)
$c.Value = "'5.268"
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -7.11%  '
$c = $ws.Cells.Item(14, 4)
$c.Value = "'6.038"
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -8.24%  '
$c = $ws.Cells.Item(15, 4)
$c.Value = "'16.44"
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -11.52%  '
$ws.Cells.Item(16, 4).Value = '1.389.54'
$ws.Cells.Item(16, 5).Value = '  -10.96%  '
$c = $ws.Cells.Item(17, 4)
$c.Value = "'0.000009867"
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -9.65%  '
$c = $ws.Cells.Item(18, 4)
$c.Value = "'0.05620"
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -14.75%  '
$c = $ws.Cells.Item(19, 4)
$c.Value = "'1.003"
$c.Style = 'Normal'
$c = $ws.Cells.Item(20, 4)
$c.Value = "'69.63"
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -16.91%  '
$c = $ws.Cells.Item(21, 4)
$c.Value = "'5.492"
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -10.25%  '
$c = $ws.Cells.Item(22, 4)
$c.Value = "'14.27"
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -7.74%  '
$c = $ws.Cells.Item(23, 4)
$c.Value = "'10.49"
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -2.21%  '
$c = $ws.Cells.Item(24, 4)
$c.Value = "'2.249"
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -4.36%  '
$ws.Cells.Item(25, 4).Value = '19.731.18'
$ws.Cells.Item(25, 5).Value = '  -8.89%  '
$c = $ws.Cells.Item(26, 4)
$c.Value = "'2.157"
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -8.52%  '
$c = $ws.Cells.Item(27, 4)
$c.Value = "'135.96"
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -9.02%  '
$c = $ws.Cells.Item(28, 4)
$c.Value = "'16.48"
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -9.99%  '
$ws.Cells.Item(29, 4).Value = '1.543.09'
$ws.Cells.Item(29, 5).Value = '  -10.90%  '
$ws.Cells.Item(30, 5).Value = '  -8.36%  '
$c = $ws.Cells.Item(31, 4)
$c.Value = "'3.825"
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -21.05%  '
$c = $ws.Cells.Item(32, 4)
$c.Value = "'5.199"
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -14.03%  '
$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.7887"
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -16.04%  '
$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.07560"
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -7.25%  '
$c = $ws.Cells.Item(35, 4)
$c.Value = "'8.170"
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -4.02%  '
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.05583"
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -7.33%  '
$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(38, 4)
$c.Value = "'4.659"
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -9.70%  '
$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.02011"
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -9.45%  '
$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.1863"
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -7.84%  '
$c = $ws.Cells.Item(41, 4)
$c.Value = "'9.913"
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -9.57%  '
$c = $ws.Cells.Item(42, 4)
$c.Value = "'1.280"
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -12.00%  '
$c = $ws.Cells.Item(43, 4)
$c.Value = "'1.044"
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -11.20%  '
$c = $ws.Cells.Item(44, 4)
$c.Value = "'3.453"
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -6.86%  '
$c = $ws.Cells.Item(45, 4)
$c.Value = "'0.5151"
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -10.60%  '
$c = $ws.Cells.Item(46, 4)
$c.Value = "'11.80"
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -9.18%  '
$c = $ws.Cells.Item(47, 4)
$c.Value = "'0.4952"
$c.Style = 'Normal'
$c = $ws.Cells.Item(48, 4)
$c.Value = "'108.14"
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -6.95%  '
$c = $ws.Cells.Item(49, 4)
$c.Value = "'1.715"
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -8.75%  '
$c = $ws.Cells.Item(50, 4)
$c.Value = "'1.002"
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +0.18%  '
$c = $ws.Cells.Item(51, 4)
$c.Value = "'1.028"
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -11.74%  '
